$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(453).Insert()

$ws.Range("A453").Value = 10
$ws.Range("B453").Value = "Vega Modelo de Temuco"
$ws.Range("C453").Value = "La Araucanía"
$ws.Range("D453").Value = 44889
$ws.Range("E453").Value = 9
$ws.Range("F453").Value = 100112023
$ws.Range("G453").Value = "Brócoli"
$ws.Range("H453").Value = "Sin especificar"
$ws.Range("I453").Value = "Primera"
$ws.Range("J453").Value = 850
$ws.Range("K453").Value = 1100
$ws.Range("L453").Value = 1100
$ws.Range("M453").Value = 1100
$ws.Range("N453").Value = "$/unidad"
$ws.Range("O453").Value = "Región del Maule"
$ws.Range("P453").Value = 1100
$ws.Range("Q453").Value = 1
$ws.Range("R453").Value = "Hortaliza"
